$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 135
$ws.Range("I13").Value = 135
$ws.Range("K13").Value = 135
$ws.Range("M13").Value = 34
$ws.Range("H28").Value = 496.38235
$ws.Range("I28").Value = 554.4138
$ws.Range("J28").Value = 159.8
$ws.Range("K28").Value = 554.4138
$ws.Range("L28").Value = 159.8
$ws.Range("M28").Value = -69.41380000000004
$ws.Range("N28").Value = -1129.8
$ws.Range("H40").Value = 19315.166
$ws.Range("I40").Value = 27250.25
$ws.Range("J40").Value = 3445
$ws.Range("K40").Value = 27250.25
$ws.Range("L40").Value = 3445
$ws.Range("M40").Value = -27075.25
$ws.Range("N40").Value = -3795
$ws.Range("H123").Value = 38849.715
$ws.Range("J123").Value = 38849.715
$ws.Range("L123").Value = 38849.715
$ws.Range("N123").Value = -48649.715
$ws.Range("H135").Value = 19231942
$ws.Range("I135").Value = 1259.35
$ws.Range("K135").Value = 11334.15
$ws.Range("M135").Value = -8799.15

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2299.8865
$ws.Range("I61").Value = 2228.1025
$ws.Range("J61").Value = 2859.8
$ws.Range("K61").Value = 2228.1025
$ws.Range("L61").Value = 2859.8
$ws.Range("M61").Value = -2016.1025
$ws.Range("N61").Value = -3283.8
$ws.Range("H133").Value = 28886.87
$ws.Range("J133").Value = 28886.87
$ws.Range("L133").Value = 28886.87
$ws.Range("N133").Value = -33946.87
$ws.Range("H136").Value = 2299.8865
$ws.Range("I136").Value = 2228.1025
$ws.Range("J136").Value = 2859.8
$ws.Range("K136").Value = 6684.3075
$ws.Range("L136").Value = 8579.400000000001
$ws.Range("M136").Value = -4134.3075
$ws.Range("N136").Value = -13679.4

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2576251.5
$ws.Range("I6").Value = 10000500
$ws.Range("J6").Value = 101502
$ws.Range("K6").Value = 10000500
$ws.Range("L6").Value = 101502
$ws.Range("M6").Value = -10000387
$ws.Range("N6").Value = -101728
$ws.Range("H58").Value = 2063.4312
$ws.Range("I58").Value = 1698.0613
$ws.Range("J58").Value = 4052.6667
$ws.Range("K58").Value = 1698.0613
$ws.Range("L58").Value = 4052.6667
$ws.Range("M58").Value = -1495.0613
$ws.Range("N58").Value = -4458.6667
$ws.Range("H122").Value = 100122940
$ws.Range("I122").Value = 166867420
$ws.Range("J122").Value = 6232
$ws.Range("K122").Value = 500602260
$ws.Range("L122").Value = 18696
$ws.Range("M122").Value = -500599810
$ws.Range("N122").Value = -23596
$ws.Range("H132").Value = 44788.305
$ws.Range("I132").Value = 1898.9474
$ws.Range("K132").Value = 5696.8422
$ws.Range("M132").Value = -3166.8422
$ws.Range("H134").Value = 2526.946
$ws.Range("I134").Value = 1058.0385
$ws.Range("J134").Value = 5998.909
$ws.Range("K134").Value = 3174.1155
$ws.Range("L134").Value = 17996.727
$ws.Range("M134").Value = -639.1155000000003
$ws.Range("N134").Value = -23066.727
$ws.Range("H136").Value = 2063.4312
$ws.Range("I136").Value = 1698.0613
$ws.Range("J136").Value = 4052.6667
$ws.Range("K136").Value = 5094.1839
$ws.Range("L136").Value = 12158.0001
$ws.Range("M136").Value = -2544.1839
$ws.Range("N136").Value = -17258.0001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H39").Value = 875.5
$ws.Range("J39").Value = 816.6667
$ws.Range("L39").Value = 2450.0001
$ws.Range("N39").Value = -3038.0001
$ws.Range("H55").Value = 1100
$ws.Range("J55").Value = 2000
$ws.Range("L55").Value = 6000
$ws.Range("N55").Value = -6354
$ws.Range("H92").Value = 1209.0625
$ws.Range("I92").Value = 1095.8572
$ws.Range("J92").Value = 2001.5
$ws.Range("K92").Value = 3287.5716
$ws.Range("L92").Value = 6004.5
$ws.Range("M92").Value = -2039.5716
$ws.Range("N92").Value = -8500.5
$ws.Range("H98").Value = 1500
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 4500
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -7496
$ws.Range("H114").Value = 1135.8889
$ws.Range("I114").Value = 974.125
$ws.Range("J114").Value = 2430
$ws.Range("K114").Value = 2922.375
$ws.Range("L114").Value = 7290
$ws.Range("M114").Value = 331.625
$ws.Range("N114").Value = -13798

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2459.4167
$ws.Range("I126").Value = 2751.2
$ws.Range("J126").Value = 2251
$ws.Range("K126").Value = 8253.599999999999
$ws.Range("L126").Value = 6753
$ws.Range("M126").Value = -5783.599999999999
$ws.Range("N126").Value = -11693

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 29115.334
$ws.Range("J123").Value = 29115.334
$ws.Range("L123").Value = 29115.334
$ws.Range("N123").Value = -38915.334
$ws.Range("H132").Value = 3699.4583
$ws.Range("I132").Value = 3585.1724
$ws.Range("J132").Value = 3873.8948
$ws.Range("K132").Value = 10755.5172
$ws.Range("L132").Value = 11621.6844
$ws.Range("M132").Value = -8225.5172
$ws.Range("N132").Value = -16681.6844
$ws.Range("H136").Value = 2771.6843
$ws.Range("I136").Value = 2406.8
$ws.Range("J136").Value = 4140
$ws.Range("K136").Value = 7220.400000000001
$ws.Range("L136").Value = 12420
$ws.Range("M136").Value = -4670.400000000001
$ws.Range("N136").Value = -17520

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 44048600
$ws.Range("I122").Value = 66072330
$ws.Range("J122").Value = 1147.5
$ws.Range("K122").Value = 198216990
$ws.Range("L122").Value = 3442.5
$ws.Range("M122").Value = -198214540
$ws.Range("N122").Value = -8342.5
$ws.Range("H132").Value = 1204.4186
$ws.Range("I132").Value = 900.5161000000001
$ws.Range("J132").Value = 1989.5
$ws.Range("K132").Value = 2701.5483
$ws.Range("L132").Value = 5968.5
$ws.Range("M132").Value = -171.5483000000004
$ws.Range("N132").Value = -11028.5
$ws.Range("H138").Value = 36687.57
$ws.Range("J138").Value = 36687.57
$ws.Range("L138").Value = 36687.57
$ws.Range("N138").Value = -46967.57
